$wb = $excel.ActiveWorkbook

# --- Rename sheets (new randomized task-order timestamps) ---
$wb.Worksheets.Item("GNG_TO-1649873095229994").Name = "GNG_TO-16502912252432835"
$wb.Worksheets.Item("NB_TO-16498730967791886").Name = "NB_TO-1650291228611419"
$wb.Worksheets.Item("RS_TO-16498730967851892").Name = "RS_TO-16502912286134226"
$wb.Worksheets.Item("TOL_TO-16498730968441892").Name = "TOL_TO-16502912286604152"
$wb.Worksheets.Item("vSAT_TO-16498730969072242").Name = "vSAT_TO-16502912287374177"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item("GNG_TO-16502912252432835")
$ws1.Range("B2").Value = "go_stims-1650291225190284.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912252102845.csv"
$ws1.Range("B4").Value = "go_stims-16502912252122889.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912252422876.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item("NB_TO-1650291228611419")
$ws2.Range("B2").Value = "ZB-match_2-16502912260002549.csv"
$ws2.Range("B3").Value = "OB-16502912264612596.csv"
$ws2.Range("B4").Value = "ZB-match_9-16502912253562834.csv"
$ws2.Range("B5").Value = "TB-16502912278591902.csv"
$ws2.Range("B6").Value = "OB-16502912267671556.csv"
$ws2.Range("B7").Value = "ZB-match_7-1650291225779256.csv"
$ws2.Range("B8").Value = "TB-16502912268541558.csv"
$ws2.Range("B9").Value = "TB-16502912285914211.csv"
$ws2.Range("B10").Value = "OB-16502912260452573.csv"

# --- Sheet 3: RS_TO ---
$ws3 = $wb.Worksheets.Item("RS_TO-16502912286134226")
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item("TOL_TO-16502912286604152")
$ws4.Range("B2").Value = "MM_stims-16502912286274276.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912286164148.csv"
$ws4.Range("B4").Value = "MM_stims-1650291228643449.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912286284146.csv"
$ws4.Range("B6").Value = "MM_stims-16502912286594145.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912286444182.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item("vSAT_TO-16502912287374177")
$ws5.Range("B2").Value = "SAT_stims-16502912286644185.csv"
$ws5.Range("B3").Value = "vSAT_stims-16502912287064483.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912287224174.csv"
$ws5.Range("B5").Value = "SAT_stims-16502912286904147.csv"
